# Excel COM-interop script: rebuild Sheet1 as a small "patients" extract
# pulled from the excel database — 2 rows x 15 columns (A:O) of numeric /
# string / date data, replacing the original 2x2 placeholder sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 1 -----------------------------------------------------------
$ws.Cells.Item(1, 1).Value  = 1
$ws.Cells.Item(1, 2).Value  = "a"
$ws.Cells.Item(1, 3).Value  = 12.33
$ws.Cells.Item(1, 4).Value  = 1.22
$ws.Cells.Item(1, 5).Value  = 1.22
$ws.Cells.Item(1, 6).Value  = 12.31
$ws.Cells.Item(1, 7).Value  = 12.43
$ws.Cells.Item(1, 8).Value  = 543.20000000000005
$ws.Cells.Item(1, 9).Value  = 123.43
$ws.Cells.Item(1, 10).Value = 54.23
$ws.Cells.Item(1, 11).Value = 134.21
$ws.Cells.Item(1, 12).Value = 1111.2
$ws.Cells.Item(1, 13).Value = 42445.895833333336
$ws.Cells.Item(1, 14).Value = 42445.895833333336
$ws.Cells.Item(1, 15).Value = 42445.895833333336

# --- Row 2 -----------------------------------------------------------
$ws.Cells.Item(2, 1).Value  = 2
$ws.Cells.Item(2, 2).Value  = "bbb"
$ws.Cells.Item(2, 3).Value  = 12.33
$ws.Cells.Item(2, 4).Value  = 1.22
$ws.Cells.Item(2, 5).Value  = 1.22
$ws.Cells.Item(2, 6).Value  = 12.31
$ws.Cells.Item(2, 7).Value  = 12.43
$ws.Cells.Item(2, 8).Value  = 543.20000000000005
$ws.Cells.Item(2, 9).Value  = 123.43
$ws.Cells.Item(2, 10).Value = 54.23
$ws.Cells.Item(2, 11).Value = 134.21
$ws.Cells.Item(2, 12).Value = 1111.2
$ws.Cells.Item(2, 13).Value = 42445.895833333336
$ws.Cells.Item(2, 14).Value = 42445.895833333336
$ws.Cells.Item(2, 15).Value = 42445.895833333336

# Date/time formatting for the M:O columns (built-in format 22 — m/d/yy h:mm)
$ws.Range("M1:O2").NumberFormat = "m/d/yy h:mm"

# Widen the date columns to fit their content
$ws.Range("M1:O2").ColumnWidth = 13.8
